$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '47.505.53'
$ws.Range("E2").Value = '  +1.29%  '

# Row 3
$ws.Range("D3").Value = '2.281.51'
$ws.Range("E3").Value = '  -2.31%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.995'
$ws.Range("E4").Value = '  -0.54%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.01'
$ws.Range("E5").Value = '  -1.60%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.63'
$ws.Range("E6").Value = '  +4.48%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.575'
$ws.Range("E7").Value = '  -0.17%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.997'
$ws.Range("E8").Value = '  -0.25%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.511'
$ws.Range("E9").Value = '  -4.35%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.30'
$ws.Range("E10").Value = '  -1.08%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0791'
$ws.Range("E11").Value = '  -2.09%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.10'
$ws.Range("E12").Value = '  -4.03%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.103'
$ws.Range("E13").Value = '  -1.32%  '

# Row 14
$ws.Range("D14").Value = '2.603.37'
$ws.Range("E14").Value = '  -3.23%  '

# Row 15
$ws.Range("D15").Value = '2.260.70'
$ws.Range("E15").Value = '  -3.29%  '

# Row 16
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.71'
$ws.Range("E16").Value = '  -2.89%  '

# Row 17
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '47.348.77'
$ws.Range("E17").Value = '  +1.32%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.804'
$ws.Range("E18").Value = '  -2.87%  '

# Row 19
$ws.Range("D19").Value = '0.0₃0977'
$ws.Range("E19").Value = '  +3.11%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.62'
$ws.Range("E20").Value = '  -7.67%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.85'
$ws.Range("E21").Value = '  -5.34%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.59'
$ws.Range("E22").Value = '  -0.90%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '249.89'
$ws.Range("E23").Value = '  +1.45%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.81'
$ws.Range("E24").Value = '  -5.64%  '

# Row 25
$ws.Range("E25").Value = '  +0.39%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.88'
$ws.Range("E26").Value = '  -5.44%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '41.42'
$ws.Range("E27").Value = '  -1.46%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.23'
$ws.Range("E28").Value = '  -3.43%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.66'
$ws.Range("E29").Value = '  -2.17%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.19'
$ws.Range("E30").Value = '  +0.41%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.83'
$ws.Range("E31").Value = '  +9.18%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.35'
$ws.Range("E32").Value = '  +5.65%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '144.80'
$ws.Range("E33").Value = '  -5.18%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.38'
$ws.Range("E34").Value = '  -6.38%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0775'
$ws.Range("E35").Value = '  -5.29%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.112'
$ws.Range("E36").Value = '  +0.61%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.116'
$ws.Range("E37").Value = '  -2.27%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '15.71'
$ws.Range("E38").Value = '  +13.92%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.68'
$ws.Range("E39").Value = '  -7.03%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.87'
$ws.Range("E40").Value = '  -3.59%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0299'
$ws.Range("E41").Value = '  -5.11%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.14'
$ws.Range("E42").Value = '  -7.59%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.993'
$ws.Range("E43").Value = '  -0.64%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '94.16'
$ws.Range("E44").Value = '  +14.32%  '

# Row 45
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.96'
$ws.Range("E45").Value = '  -1.41%  '

# Row 46
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '1.793.24'
$ws.Range("E46").Value = '  -2.37%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '71.44'
$ws.Range("E47").Value = '  -3.85%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.186'
$ws.Range("E48").Value = '  -5.46%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.84'
$ws.Range("E49").Value = '  -2.23%  '

# Row 50
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.92'
$ws.Range("E50").Value = '  -1.45%  '

# Row 51
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '94.77'
$ws.Range("E51").Value = '  -3.96%  '

